# Refresh the NATMI Col18a1-Itgb1 LR-pair table (natmiOut/YoungD0/LR-pairs_lrc2p)
# following Dr Hou's advice: ECs are now also included as a sending cluster,
# and edge expression weights/specificities are recomputed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one data row: Sending cluster, Ligand symbol, Receptor
# symbol, Target cluster, then the 16 numeric metric columns (E:T).
$rows = @(
    @("ECs", "Col18a1", "Itgb1", "ECs", 2, 0.6666666666666666, 10.92859066666667, 32.785772, 0.2185558471001832, 0.2185558471001832, 3, 1, 112.513392, 337.540176, 0.3275312977368564, 0.3275312977368564, 1229.612805686208, 11066.51525117587, 0.07158388022870095, 0.07158388022870095),
    @("ECs", "Col18a1", "Itgb1", "FAPs", 2, 0.6666666666666666, 10.92859066666667, 32.785772, 0.2185558471001832, 0.2185558471001832, 3, 1, 106.314466, 318.943398, 0.3094859589441663, 0.3094859589441664, 1161.867280859251, 10456.80552773326, 0.06763996592265478, 0.0676399659226548),
    @("ECs", "Col18a1", "Itgb1", "sCs", 2, 0.6666666666666666, 10.92859066666667, 32.785772, 0.2185558471001832, 0.2185558471001832, 3, 1, 124.6916553333333, 374.074966, 0.3629827433189773, 0.3629827433189773, 1362.704060687084, 12264.33654618375, 0.07933200094882743, 0.07933200094882746),
    @("FAPs", "Col18a1", "Itgb1", "ECs", 3, 1, 21.275312, 63.825936, 0.4254751576214852, 0.4254751576214852, 3, 1, 112.513392, 337.540176, 0.3275312977368564, 0.3275312977368564, 2393.757518978304, 21543.81767080473, 0.1393564305305585, 0.1393564305305586),
    @("FAPs", "Col18a1", "Itgb1", "FAPs", 3, 1, 21.275312, 63.825936, 0.4254751576214852, 0.4254751576214852, 3, 1, 106.314466, 318.943398, 0.3094859589441663, 0.3094859589441664, 2261.873434263392, 20356.86090837053, 0.1316785871634057, 0.1316785871634057),
    @("FAPs", "Col18a1", "Itgb1", "sCs", 3, 1, 21.275312, 63.825936, 0.4254751576214852, 0.4254751576214852, 3, 1, 124.6916553333333, 374.074966, 0.3629827433189773, 0.3629827433189773, 2652.853871013131, 23875.68483911818, 0.154440139927521, 0.154440139927521),
    @("sCs", "Col18a1", "Itgb1", "ECs", 3, 1, 17.79975, 53.39925, 0.3559689952783316, 0.3559689952783316, 3, 1, 112.513392, 337.540176, 0.3275312977368564, 0.3275312977368564, 2002.710249252, 18024.392243268, 0.1165909869775968, 0.1165909869775969),
    @("sCs", "Col18a1", "Itgb1", "FAPs", 3, 1, 17.79975, 53.39925, 0.3559689952783316, 0.3559689952783316, 3, 1, 106.314466, 318.943398, 0.3094859589441663, 0.3094859589441664, 1892.3709161835, 17031.3382456515, 0.1101674058581059, 0.1101674058581059),
    @("sCs", "Col18a1", "Itgb1", "sCs", 3, 1, 17.79975, 53.39925, 0.3559689952783316, 0.3559689952783316, 3, 1, 124.6916553333333, 374.074966, 0.3629827433189773, 0.3629827433189773, 2219.4802920195, 19975.3226281755, 0.1292106024426289, 0.1292106024426289)
)

$numRows = $rows.Count
$numCols = $rows[0].Count

$data = New-Object 'object[,]' $numRows,$numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $data[$r, $c] = $rows[$r][$c]
    }
}

# The table now spans 9 data rows (was 6) -> worksheet dimension becomes A1:T10.
$ws.Range("A2:T10").Value = $data

